$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.06370000000000001
$ws.Range("E2").Value = 1.353
$ws.Range("G2").Value = 0.279773156899811
$ws.Range("H2").Value = 0.279773156899811
$ws.Range("I2").Value = 0.2879962192816635
$ws.Range("J2").Value = 0.167709459143527
$ws.Range("K2").Value = 1.905
$ws.Range("L2").Value = 0.1800567107750472
$ws.Range("M2").Value = 2.484
$ws.Range("N2").Value = 0.08619014573213046
$ws.Range("O2").Value = 1.303937007874016
$ws.Range("P2").Value = 2.064
$ws.Range("Q2").Value = 0.07161693268563497
$ws.Range("R2").Value = 1.083464566929134
$ws.Range("S2").Value = 0.4200000000000002
$ws.Range("T2").Value = 0.1690821256038648
$ws.Range("U2").Value = 4.69
$ws.Range("V2").Value = 0.162734212352533
$ws.Range("W2").Value = 0.186392183931554
$ws.Range("X2").Value = 0.06114138249983702
$ws.Range("Y2").Value = 0.125250801431717
$ws.Range("Z2").Value = 0.727147766323024
$ws.Range("AA2").Value = 0.1755335665699866
$ws.Range("AB2").Value = 0.0572041771518761
$ws.Range("AC2").Value = 0.1183293894181105
$ws.Range("AD2").Value = 3.715
$ws.Range("AF2").Value = 3.715
$ws.Range("AG2").Value = -0.9750000000000005
$ws.Range("AH2").Value = 0.1141847241432304
$ws.Range("AI2").Value = 0.2104219767771169
$ws.Range("AJ2").Value = -0.03501526306338663
$ws.Range("AK2").Value = -0.07520246818357119
$ws.Range("AL2").Value = 0.121
$ws.Range("AM2").Value = 0.115
$ws.Range("AN2").Value = 1.255067567567568
$ws.Range("AO2").Value = 25.18181818181818
$ws.Range("AP2").Value = -0.3293918918918921
$ws.Range("AQ2").Value = 26.49565217391305

# Row 3
$ws.Range("D3").Value = 0.06370000000000001
$ws.Range("E3").Value = 1.353
$ws.Range("G3").Value = 0.3234972677595628
$ws.Range("H3").Value = 0.3234972677595628
$ws.Range("I3").Value = 0.3136612021857924
$ws.Range("J3").Value = 0.2084794068238429
$ws.Range("K3").Value = 1.88
$ws.Range("L3").Value = 0.2054644808743169
$ws.Range("M3").Value = 2.116
$ws.Range("N3").Value = 0.1022222222222222
$ws.Range("O3").Value = 1.125531914893617
$ws.Range("P3").Value = 1.89
$ws.Range("Q3").Value = 0.09130434782608696
$ws.Range("R3").Value = 1.00531914893617
$ws.Range("S3").Value = 0.2260000000000002
$ws.Range("T3").Value = 0.1068052930056712
$ws.Range("U3").Value = 2.47
$ws.Range("V3").Value = 0.1193236714975845
$ws.Range("W3").Value = 0.3700787401574803
$ws.Range("X3").Value = 0.05584132398743984
$ws.Range("Y3").Value = 0.3142374161700404
$ws.Range("Z3").Value = 1.636558755142193
$ws.Range("AA3").Value = 0.3411887985044111
$ws.Range("AB3").Value = 0.0551564865285958
$ws.Range("AC3").Value = 0.2860323119758154
$ws.Range("AD3").Value = 0.865
$ws.Range("AF3").Value = 0.865
$ws.Range("AG3").Value = -1.605
$ws.Range("AH3").Value = 0.04011129144447021
$ws.Range("AI3").Value = 0.1409942950285249
$ws.Range("AJ3").Value = -0.08405341712490182
$ws.Range("AK3").Value = -0.4379263301500684
$ws.Range("AL3").Value = 0.043
$ws.Range("AM3").Value = 0.037
$ws.Range("AN3").Value = 0.2922297297297297
$ws.Range("AO3").Value = 66.74418604651163
$ws.Range("AP3").Value = -0.5422297297297298
$ws.Range("AQ3").Value = 77.56756756756758

# Row 4
$ws.Range("I4").Value = 0.1237762237762238
$ws.Range("J4").Value = 0.06188811188811189
$ws.Range("K4").Value = 0.025
$ws.Range("L4").Value = 0.01748251748251748
$ws.Range("M4").Value = 0.368
$ws.Range("N4").Value = 0.04532019704433498
$ws.Range("O4").Value = 14.72
$ws.Range("P4").Value = 0.174
$ws.Range("Q4").Value = 0.02142857142857143
$ws.Range("R4").Value = 6.959999999999999
$ws.Range("S4").Value = 0.194
$ws.Range("T4").Value = 0.5271739130434783
$ws.Range("U4").Value = 2.22
$ws.Range("V4").Value = 0.2733990147783252
$ws.Range("W4").Value = 0.002705627705627706
$ws.Range("X4").Value = 0.0664414410122342
$ws.Range("Y4").Value = -0.06373581330660649
$ws.Range("Z4").Value = 0.1596160285746177
$ws.Range("AA4").Value = 0.009878334635562005
$ws.Range("AB4").Value = 0.05925186777515639
$ws.Range("AC4").Value = -0.04937353313959439
$ws.Range("AD4").Value = 2.85
$ws.Range("AF4").Value = 2.85
$ws.Range("AG4").Value = 0.6299999999999999
$ws.Range("AH4").Value = 0.2597994530537831
$ws.Range("AI4").Value = 0.2473958333333333
$ws.Range("AJ4").Value = 0.07199999999999999
$ws.Range("AK4").Value = 0.06774193548387095
$ws.Range("AL4").Value = 0.078
$ws.Range("AM4").Value = 0.078
$ws.Range("AO4").Value = 2.269230769230769
$ws.Range("AQ4").Value = 2.269230769230769
